# Update Sheets via scheduled runner: apply recalculated profit figures
$wb = $excel.ActiveWorkbook

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1067.8864
$ws.Range("I129").Value = 362.33334
$ws.Range("J129").Value = 1092.7882
$ws.Range("K129").Value = 1087.00002
$ws.Range("L129").Value = 3278.3646
$ws.Range("M129").Value = 3912.99998
$ws.Range("N129").Value = -13278.3646

# ALC row 133
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 69558.336
$ws.Range("J133").Value = 69558.336
$ws.Range("L133").Value = 69558.336
$ws.Range("N133").Value = -79678.336

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2127.2046
$ws.Range("I137").Value = 1361.6666
$ws.Range("J137").Value = 2826.1738
$ws.Range("K137").Value = 4084.9998
$ws.Range("L137").Value = 8478.5214
$ws.Range("M137").Value = -1534.9998
$ws.Range("N137").Value = -13578.5214

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2514.0613
$ws.Range("I138").Value = 1255.4634
$ws.Range("J138").Value = 3419.3684
$ws.Range("K138").Value = 3766.3902
$ws.Range("L138").Value = 10258.1052
$ws.Range("M138").Value = 1373.6098
$ws.Range("N138").Value = -20538.1052

# ARM row 24
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 42785
$ws.Range("J24").Value = 42785
$ws.Range("L24").Value = 42785
$ws.Range("N24").Value = -43533

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9859.219999999999
$ws.Range("I32").Value = 10236.528
$ws.Range("J32").Value = 7658.25
$ws.Range("K32").Value = 10236.528
$ws.Range("L32").Value = 7658.25
$ws.Range("M32").Value = -9949.528
$ws.Range("N32").Value = -8232.25

# ARM row 100
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H100").Value = 42785
$ws.Range("J100").Value = 42785
$ws.Range("L100").Value = 42785
$ws.Range("N100").Value = -44949

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3437
$ws.Range("I105").Value = 3332.8333
$ws.Range("J105").Value = 3749.5
$ws.Range("K105").Value = 3332.8333
$ws.Range("L105").Value = 3749.5
$ws.Range("M105").Value = -1585.8333
$ws.Range("N105").Value = -7243.5

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2666.1562
$ws.Range("I134").Value = 2305.75
$ws.Range("J134").Value = 3747.375
$ws.Range("K134").Value = 6917.25
$ws.Range("L134").Value = 11242.125
$ws.Range("M134").Value = -4382.25
$ws.Range("N134").Value = -16312.125

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1808.5438
$ws.Range("I31").Value = 2244.5715
$ws.Range("J31").Value = 1387.5518
$ws.Range("K31").Value = 2244.5715
$ws.Range("L31").Value = 1387.5518
$ws.Range("M31").Value = -1949.5715
$ws.Range("N31").Value = -1977.5518

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1808.5438
$ws.Range("I34").Value = 2244.5715
$ws.Range("J34").Value = 1387.5518
$ws.Range("K34").Value = 2244.5715
$ws.Range("L34").Value = 1387.5518
$ws.Range("M34").Value = -2042.5715
$ws.Range("N34").Value = -1791.5518

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1827.2222
$ws.Range("I99").Value = 1864.6154
$ws.Range("J99").Value = 1730
$ws.Range("K99").Value = 1864.6154
$ws.Range("L99").Value = 1730
$ws.Range("M99").Value = -366.6153999999999
$ws.Range("N99").Value = -4726

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1827.2222
$ws.Range("I126").Value = 1864.6154
$ws.Range("J126").Value = 1730
$ws.Range("K126").Value = 5593.8462
$ws.Range("L126").Value = 5190
$ws.Range("M126").Value = -3123.8462
$ws.Range("N126").Value = -10130

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2256903.8
$ws.Range("I132").Value = 6759712.5
$ws.Range("J132").Value = 5499.5
$ws.Range("K132").Value = 20279137.5
$ws.Range("L132").Value = 16498.5
$ws.Range("M132").Value = -20276607.5
$ws.Range("N132").Value = -21558.5

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 865.89
$ws.Range("I68").Value = 614.4761999999999
$ws.Range("J68").Value = 1293.973
$ws.Range("K68").Value = 1843.4286
$ws.Range("L68").Value = 3881.919
$ws.Range("M68").Value = -1032.4286
$ws.Range("N68").Value = -5503.919

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 865.89
$ws.Range("I71").Value = 614.4761999999999
$ws.Range("J71").Value = 1293.973
$ws.Range("K71").Value = 5530.2858
$ws.Range("L71").Value = 11645.757
$ws.Range("M71").Value = -1474.2858
$ws.Range("N71").Value = -19757.757

# CUL row 86
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1828.7778
$ws.Range("I86").Value = 2137
$ws.Range("J86").Value = 750
$ws.Range("K86").Value = 6411
$ws.Range("L86").Value = 2250
$ws.Range("M86").Value = -5225
$ws.Range("N86").Value = -4622

# CUL row 89
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 1828.7778
$ws.Range("I89").Value = 2137
$ws.Range("J89").Value = 750
$ws.Range("K89").Value = 19233
$ws.Range("L89").Value = 6750
$ws.Range("M89").Value = -13305
$ws.Range("N89").Value = -18606

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 769.5700000000001
$ws.Range("I107").Value = 654.5769
$ws.Range("J107").Value = 809.9729599999999
$ws.Range("K107").Value = 1963.7307
$ws.Range("L107").Value = 2429.91888
$ws.Range("M107").Value = -43.73070000000007
$ws.Range("N107").Value = -6269.918879999999

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 439.7
$ws.Range("I122").Value = 574
$ws.Range("J122").Value = 337
$ws.Range("K122").Value = 5166
$ws.Range("L122").Value = 3033
$ws.Range("M122").Value = -2716
$ws.Range("N122").Value = -7933

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 21279816
$ws.Range("I131").Value = 577.7273
$ws.Range("J131").Value = 27781804
$ws.Range("K131").Value = 1733.1819
$ws.Range("L131").Value = 83345412
$ws.Range("M131").Value = 3306.8181
$ws.Range("N131").Value = -83355492

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1021.7895
$ws.Range("I132").Value = 937.55554
$ws.Range("J132").Value = 1047.931
$ws.Range("K132").Value = 8437.99986
$ws.Range("L132").Value = 9431.379000000001
$ws.Range("M132").Value = -5907.99986
$ws.Range("N132").Value = -14491.379

# GSM row 14
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 19200200
$ws.Range("J14").Value = 1000
$ws.Range("L14").Value = 1000
$ws.Range("N14").Value = -1336

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2859.4285
$ws.Range("I113").Value = 1750
$ws.Range("J113").Value = 3303.2
$ws.Range("K113").Value = 1750
$ws.Range("L113").Value = 3303.2
$ws.Range("M113").Value = 420
$ws.Range("N113").Value = -7643.2

# GSM row 123
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 8618.379000000001
$ws.Range("J123").Value = 8618.379000000001
$ws.Range("L123").Value = 8618.379000000001
$ws.Range("N123").Value = -13518.379

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3300
$ws.Range("I40").Value = 3300
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3300
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3164
$ws.Range("N40").ClearContents()

# LTW row 115
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H115").Value = 89514.86
$ws.Range("J115").Value = 89514.86
$ws.Range("L115").Value = 89514.86
$ws.Range("N115").Value = -91864.86

# LTW row 120
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H120").Value = 24698.715
$ws.Range("J120").Value = 24698.715
$ws.Range("L120").Value = 24698.715
$ws.Range("N120").Value = -34374.715

# WVR row 64
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 29923.076
$ws.Range("J64").Value = 29923.076
$ws.Range("L64").Value = 29923.076
$ws.Range("N64").Value = -30419.076

# WVR row 67
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 29923.076
$ws.Range("J67").Value = 29923.076
$ws.Range("L67").Value = 29923.076
$ws.Range("N67").Value = -31639.076

# WVR row 120
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H120").Value = 38000
$ws.Range("J120").Value = 38000
$ws.Range("L120").Value = 38000
$ws.Range("N120").Value = -47676

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1749.9149
$ws.Range("I132").Value = 1537.5834
$ws.Range("J132").Value = 2444.818
$ws.Range("K132").Value = 4612.7502
$ws.Range("L132").Value = 7334.454000000001
$ws.Range("M132").Value = -2082.7502
$ws.Range("N132").Value = -12394.454
